# Update countries & provincias Spain
# Refresh the COVID-19 country statistics table:
#  - bump the "last updated" timestamp
#  - update case counts for several countries (new data pull)
#  - re-insert/re-rank Serbia, Bielorrusia and Sri Lanka, which causes the
#    rows between their old and new rank to shift by one position
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Abril de 2020 a las 15:22"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 368174
$ws.Cells.Item(4, 3).Value = 1170
$ws.Cells.Item(4, 4).Value = 19814
$ws.Cells.Item(4, 5).Value = 337394
$ws.Cells.Item(4, 6).Value = 8983
$ws.Cells.Item(4, 7).Value = 95
$ws.Cells.Item(4, 8).Value = 10966

# Row 17
$ws.Cells.Item(17, 1).Value = "Austria"
$ws.Cells.Item(17, 2).Value = 12519
$ws.Cells.Item(17, 3).Value = 222
$ws.Cells.Item(17, 4).Value = 4046
$ws.Cells.Item(17, 5).Value = 8230
$ws.Cells.Item(17, 6).Value = 243
$ws.Cells.Item(17, 7).Value = 23
$ws.Cells.Item(17, 8).Value = 243

# Row 28
$ws.Cells.Item(28, 1).Value = "India"
$ws.Cells.Item(28, 2).Value = 4911
$ws.Cells.Item(28, 3).Value = 133
$ws.Cells.Item(28, 4).Value = 382
$ws.Cells.Item(28, 5).Value = 4392
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 1
$ws.Cells.Item(28, 8).Value = 137

# Row 42
$ws.Cells.Item(42, 1).Value = "Serbia"
$ws.Cells.Item(42, 2).Value = 2447
$ws.Cells.Item(42, 3).Value = 247
$ws.Cells.Item(42, 4).Value = 118
$ws.Cells.Item(42, 5).Value = 2268
$ws.Cells.Item(42, 6).Value = 109
$ws.Cells.Item(42, 7).Value = 3
$ws.Cells.Item(42, 8).Value = 61

# Row 43
$ws.Cells.Item(43, 1).Value = "Mexico"
$ws.Cells.Item(43, 2).Value = 2439
$ws.Cells.Item(43, 3).Value = 296
$ws.Cells.Item(43, 4).Value = 633
$ws.Cells.Item(43, 5).Value = 1681
$ws.Cells.Item(43, 6).Value = 89
$ws.Cells.Item(43, 7).Value = 31
$ws.Cells.Item(43, 8).Value = 125

# Row 44
$ws.Cells.Item(44, 1).Value = "Finlandia"
$ws.Cells.Item(44, 2).Value = 2308
$ws.Cells.Item(44, 3).Value = 132
$ws.Cells.Item(44, 4).Value = 300
$ws.Cells.Item(44, 5).Value = 1974
$ws.Cells.Item(44, 6).Value = 81
$ws.Cells.Item(44, 7).Value = 7
$ws.Cells.Item(44, 8).Value = 34

# Row 45
$ws.Cells.Item(45, 1).Value = "Tailandia"
$ws.Cells.Item(45, 2).Value = 2258
$ws.Cells.Item(45, 3).Value = 38
$ws.Cells.Item(45, 4).Value = 824
$ws.Cells.Item(45, 5).Value = 1407
$ws.Cells.Item(45, 6).Value = 30
$ws.Cells.Item(45, 7).Value = 1
$ws.Cells.Item(45, 8).Value = 27

# Row 52
$ws.Cells.Item(52, 1).Value = "Argentina"
$ws.Cells.Item(52, 2).Value = 1628
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 338
$ws.Cells.Item(52, 5).Value = 1236
$ws.Cells.Item(52, 6).Value = 96
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 54

# Row 68
$ws.Cells.Item(68, 1).Value = "Bielorrusia"
$ws.Cells.Item(68, 2).Value = 860
$ws.Cells.Item(68, 3).Value = 160
$ws.Cells.Item(68, 4).Value = 53
$ws.Cells.Item(68, 5).Value = 794
$ws.Cells.Item(68, 6).Value = 31
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 13

# Row 69
$ws.Cells.Item(69, 1).Value = "Armenia"
$ws.Cells.Item(69, 2).Value = 853
$ws.Cells.Item(69, 3).Value = 20
$ws.Cells.Item(69, 4).Value = 87
$ws.Cells.Item(69, 5).Value = 758
$ws.Cells.Item(69, 6).Value = 30
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 8

# Row 70
$ws.Cells.Item(70, 1).Value = "Hungria"
$ws.Cells.Item(70, 2).Value = 817
$ws.Cells.Item(70, 3).Value = 73
$ws.Cells.Item(70, 4).Value = 71
$ws.Cells.Item(70, 5).Value = 699
$ws.Cells.Item(70, 6).Value = 17
$ws.Cells.Item(70, 7).Value = 9
$ws.Cells.Item(70, 8).Value = 47

# Row 71
$ws.Cells.Item(71, 1).Value = "Barein"
$ws.Cells.Item(71, 2).Value = 756
$ws.Cells.Item(71, 3).Value = 0
$ws.Cells.Item(71, 4).Value = 458
$ws.Cells.Item(71, 5).Value = 294
$ws.Cells.Item(71, 6).Value = 4
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 4

# Row 72
$ws.Cells.Item(72, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(72, 2).Value = 754
$ws.Cells.Item(72, 3).Value = 80
$ws.Cells.Item(72, 4).Value = 68
$ws.Cells.Item(72, 5).Value = 653
$ws.Cells.Item(72, 6).Value = 4
$ws.Cells.Item(72, 7).Value = 4
$ws.Cells.Item(72, 8).Value = 33

# Row 73
$ws.Cells.Item(73, 1).Value = "Kuwait"
$ws.Cells.Item(73, 2).Value = 743
$ws.Cells.Item(73, 3).Value = 78
$ws.Cells.Item(73, 4).Value = 105
$ws.Cells.Item(73, 5).Value = 637
$ws.Cells.Item(73, 6).Value = 23
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 1

# Row 74
$ws.Cells.Item(74, 1).Value = "Azerbaiyan"
$ws.Cells.Item(74, 2).Value = 717
$ws.Cells.Item(74, 3).Value = 76
$ws.Cells.Item(74, 4).Value = 44
$ws.Cells.Item(74, 5).Value = 665
$ws.Cells.Item(74, 6).Value = 23
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 8

# Row 75
$ws.Cells.Item(75, 1).Value = "Crucero"
$ws.Cells.Item(75, 2).Value = 712
$ws.Cells.Item(75, 3).Value = 0
$ws.Cells.Item(75, 4).Value = 619
$ws.Cells.Item(75, 5).Value = 82
$ws.Cells.Item(75, 6).Value = 10
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 11

# Row 112
$ws.Cells.Item(112, 1).Value = "Sri Lanka"
$ws.Cells.Item(112, 2).Value = 185
$ws.Cells.Item(112, 3).Value = 7
$ws.Cells.Item(112, 4).Value = 42
$ws.Cells.Item(112, 5).Value = 137
$ws.Cells.Item(112, 6).Value = 5
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = 6

# Row 113
$ws.Cells.Item(113, 1).Value = "Islas Feroe"
$ws.Cells.Item(113, 2).Value = 184
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(113, 4).Value = 129
$ws.Cells.Item(113, 5).Value = 55
$ws.Cells.Item(113, 6).Value = 2
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 0

# Row 114
$ws.Cells.Item(114, 1).Value = "Kenia"
$ws.Cells.Item(114, 2).Value = 172
$ws.Cells.Item(114, 3).Value = 14
$ws.Cells.Item(114, 4).Value = 7
$ws.Cells.Item(114, 5).Value = 159
$ws.Cells.Item(114, 6).Value = 2
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 6
